# "Integrate student availability and assign them in group"
#
# The sheet originally held 10 students (rows 2-11): 5 "BCS" students
# (rows 2-6) and 5 "BSE" students (rows 7-11). This edit replaces the
# 5 existing BCS placeholder rows with a full cohort of 37 BCS students
# (newly available StudentIDs, counting down from 24000433 in steps of
# 5), pushing the original 5 BSE rows down to the bottom of the table
# (rows 39-43) unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 5 original BSE rows (StudentID, Programme) - preserved as-is, just
# relocated further down the sheet.
$bseStudents = @(24000004, 24000009, 24000014, 24000019, 24000024)

# The newly integrated BCS cohort: 37 students, StudentID counting down
# from 24000433 in steps of 5.
$bcsStudents = @()
$nextId = 24000433
for ($i = 0; $i -lt 37; $i++) {
    $bcsStudents += $nextId
    $nextId -= 5
}

$row = 2

foreach ($studentId in $bcsStudents) {
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = $studentId
    $ws.Cells.Item($row, 3).Value = "BCS"
    $row++
}

foreach ($studentId in $bseStudents) {
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = $studentId
    $ws.Cells.Item($row, 3).Value = "BSE"
    $row++
}

# Reflect the resulting scroll/zoom/selection state of the sheet view.
$excel.ActiveWindow.Zoom = 130
$ws.Range("C35").Select()
